# Update "paises.xlsx" (sheet "Pais") per COVID data refresh + 2 ranking swaps.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Julio de 2020 a las 03:23"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 3040833
$ws.Range("C4").Value = 50586
$ws.Range("D4").Value = 1324947
$ws.Range("E4").Value = 1582907
$ws.Range("G4").Value = 378
$ws.Range("H4").Value = 132979

# --- Row 23: Canada ---
$ws.Range("B23").Value = 105935
$ws.Range("C23").Value = 399
$ws.Range("E23").Value = 27672

# --- Row 38: Kazajistan ---
$ws.Range("E38").Value = 20976
$ws.Range("G38").Value = 13
$ws.Range("H38").Value = 264

# --- Row 73: Noruega ---
$ws.Range("B73").Value = 8936
$ws.Range("C73").Value = 6
$ws.Range("E73").Value = 547

# --- Rows 81/82: Venezuela overtakes Finlandia in ranking ---
# Row 81 was Finlandia, now becomes Venezuela with its updated totals.
$ws.Range("A81").Value = "Venezuela"
$ws.Range("B81").Value = 7411
$ws.Range("C81").Value = 242
$ws.Range("D81").Value = 2100
$ws.Range("E81").Value = 5243
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = 68

# Row 82 was Venezuela, now becomes Finlandia with its (unchanged) prior totals.
$ws.Range("A82").Value = "Finlandia"
$ws.Range("B82").Value = 7257
$ws.Range("C82").Value = 4
$ws.Range("D82").Value = 6700
$ws.Range("E82").Value = 228
$ws.Range("H82").Value = 329

# --- Row 126: Cabo Verde ---
$ws.Range("D126").Value = 722
$ws.Range("E126").Value = 724

# --- Rows 209/210: Islas Malvinas / Groenlandia swap (tied totals, order flips) ---
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"
